$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell A1 ref, new literal text value.
# Values are written with a leading apostrophe (classic Excel "force text"
# entry prefix) so numeric-looking strings (e.g. '607.61') are NOT
# auto-coerced to numbers, matching the source data which stores every
# Price/Volume column cell as literal text. Style is reset to "Normal"
# right after so no stray NumberFormat/quotePrefix styling is introduced.
$edits = @(
    @('D2', '66.398.23'),
    @('E2', '  +0.35%  '),
    @('D3', '3.595.86'),
    @('E3', '  +0.98%  '),
    @('E4', '  -0.04%  '),
    @('D5', '607.61'),
    @('E5', '  +0.30%  '),
    @('D6', '148.37'),
    @('E6', '  +2.77%  '),
    @('D7', '3.595.53'),
    @('E7', '  +1.05%  '),
    @('E8', '  -0.06%  '),
    @('D9', '0.487'),
    @('E9', '  -0.66%  '),
    @('E10', '  +0.11%  '),
    @('D11', '7.87'),
    @('E11', '  +0.44%  '),
    @('E12', '  +0.59%  '),
    @('D13', '4.203.86'),
    @('E13', '  +0.92%  '),
    @('E14', '  -0.41%  '),
    @('D15', '29.57'),
    @('E15', '  -1.43%  '),
    @('D16', '3.592.06'),
    @('E16', '  +0.42%  '),
    @('E17', '  +2.09%  '),
    @('D18', '66.440.74'),
    @('E18', '  +0.30%  '),
    @('D19', '11.11'),
    @('E19', '  -2.97%  '),
    @('E20', '  +2.06%  '),
    @('D21', '14.89'),
    @('E21', '  +1.27%  '),
    @('D22', '423.01'),
    @('E22', '  -1.68%  '),
    @('E23', '  +0.43%  '),
    @('D24', '78.68'),
    @('E24', '  -1.47%  '),
    @('E25', '  +0.10%  '),
    @('D26', '0.0000121'),
    @('E26', '  +3.62%  '),
    @('D27', '8.24'),
    @('E27', '  +5.00%  '),
    @('D28', '9.39'),
    @('E28', '  +3.00%  '),
    @('D29', '2.50'),
    @('E29', '  +0.06%  '),
    @('D30', '1.00'),
    @('E30', '  +0.03%  '),
    @('D31', '3.592.33'),
    @('E31', '  +0.90%  '),
    @('D32', '0.158'),
    @('E32', '  +3.91%  '),
    @('E33', '  -0.20%  '),
    @('E34', '  -1.23%  '),
    @('E35', '  +0.00%  '),
    @('E36', '  -0.39%  '),
    @('D37', '5.59'),
    @('E37', '  +0.93%  '),
    @('D39', '175.05'),
    @('E39', '  +1.11%  '),
    @('D40', '0.0852'),
    @('E40', '  +0.85%  '),
    @('D41', '5.19'),
    @('E41', '  +0.38%  '),
    @('E42', '  -0.79%  '),
    @('D43', '45.95'),
    @('E43', '  -0.17%  '),
    @('E44', '  -4.39%  '),
    @('D45', '0.999'),
    @('E45', '  -0.04%  '),
    @('D46', '2.52'),
    @('E46', '  +4.79%  '),
    @('D47', '23.81'),
    @('E47', '  +3.59%  '),
    @('B48', 'InjectiveProtocol'),
    @('C48', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @('D48', '24.26'),
    @('E48', '  -2.80%  '),
    @('B49', 'Cosmos'),
    @('C49', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @('D49', '7.15'),
    @('E49', '  +0.42%  '),
    @('E50', '  -5.59%  '),
    @('D51', '0.957'),
    @('E51', '  +2.97%  ')
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit[0])
    $cell.Value = "'" + $edit[1]
    $cell.Style = "Normal"
}
